$d = $word.ActiveDocument

# ============================================================
# Change 1: merge the split "*Índice por ..." runs in the
# first ("Linear Probing") bullet list.
# Each of these paragraphs only has the two runs being merged,
# so a straight literal Find/Replace is safe (nothing after the
# match inside the paragraph gets swallowed).
# ============================================================
$d.Content.Find.Execute(
    "*Índice por nombre del artista con una lista con sus obras como valor",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "*Índice por nombre del artista con una lista con sus obras como valor",
    2) | Out-Null

$d.Content.Find.Execute(
    "*Índice por nacionalidades que tiene como valor la lista de obras",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "*Índice por nacionalidades que tiene como valor la lista de obras",
    2) | Out-Null

$d.Content.Find.Execute(
    "*Índice por departamentos del museo",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "*Índice por departamentos del museo",
    2) | Out-Null

# ============================================================
# Change 2: append a new run containing a single space right
# after "...usaremos Separate Chaining." — nothing before that
# point (including the spell-check proofErr markers around
# "Separate"/"Chaining") may be touched, so we use InsertAfter
# on a collapsed range instead of a Find/Replace (which would
# otherwise coalesce every run it passes through).
# ============================================================
$rng = $d.Content
$rng.Find.Execute(
    "también usaremos Separate Chaining.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$insPoint = $d.Range($rng.End, $rng.End)
$insPoint.InsertAfter(" ")

# ============================================================
# Change 3: merge the split "*Índice por id de los artistas..."
# run in the "Linear Probing" list. Whole paragraph content, so
# a literal Find/Replace is safe.
# ============================================================
$d.Content.Find.Execute(
    "*Índice por id de los artistas con sus diccionarios como valores",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "*Índice por id de los artistas con sus diccionarios como valores",
    2) | Out-Null

# ============================================================
# Change 4: merge the "*Índice ..." + description runs in the
# load-factor list.
# ============================================================

# 4a: "nombre del artista" bullet - whole paragraph, safe as-is.
$d.Content.Find.Execute(
    "*Índice por nombre del artista con una lista con sus obras como valor: el número de obras de un artista sobre el número de artistas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "*Índice por nombre del artista con una lista con sus obras como valor: el número de obras de un artista sobre el número de artistas.",
    2) | Out-Null

# 4b: "nacionalidades" bullet - whole paragraph, safe as-is.
$d.Content.Find.Execute(
    "*Índice por nacionalidades que tiene como valor la lista de obras: el número de obras que tienen algún artista de una nacionalidad específica sobre el número de nacionalidades en el catálogo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "*Índice por nacionalidades que tiene como valor la lista de obras: el número de obras que tienen algún artista de una nacionalidad específica sobre el número de nacionalidades en el catálogo.",
    2) | Out-Null

# 4c: "departamentos del museo" bullet - the paragraph has a
# THIRD run ("departamentos.") after the text being merged, and
# this engine's replace coalesces every run it crosses *plus*
# any further run-siblings up to the next non-<w:r> element (or
# end of paragraph). To stop it from swallowing the trailing
# "departamentos." run we drop a temporary bookmark right before
# it (bookmarkStart/End act as a hard boundary) and remove the
# bookmark again once the merge is done.
$rng2 = $d.Content
$rng2.Find.Execute(
    "total de departamentos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$depStart = $rng2.End - [string]"departamentos.".Length
$d.Bookmarks.Add("zzblock1", $d.Range($depStart, $depStart)) | Out-Null

$d.Content.Find.Execute(
    "*Índice por departamentos del museo: el número de obras para un departamento sobre el número total de ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "*Índice por departamentos del museo: el número de obras para un departamento sobre el número total de ",
    2) | Out-Null

$d.Bookmarks("zzblock1").Delete()

# ============================================================
# Change 5: merge "*Índice por id de los artistas..." + ": "
# runs in the load-factor list. Same situation as 4c — a third
# run ("en este caso, ...") follows and must stay untouched, so
# we use the same bookmark-blocker trick.
# ============================================================
$rng3 = $d.Content
$rng3.Find.Execute(
    "valores: en este caso",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$enEsteStart = $rng3.End - [string]"en este caso".Length
$d.Bookmarks.Add("zzblock2", $d.Range($enEsteStart, $enEsteStart)) | Out-Null

$d.Content.Find.Execute(
    "*Índice por id de los artistas con sus diccionarios como valores: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "*Índice por id de los artistas con sus diccionarios como valores: ",
    2) | Out-Null

$d.Bookmarks("zzblock2").Delete()
